$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.147.63"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.880.71"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'482.14"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'144.60"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.722"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("E10").Value = "  +4.43%  "
$ws.Range("D11").Value = "'0.0000352"
$ws.Range("E11").Value = "  +12.01%  "
$ws.Range("D12").Value = "'42.63"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "'10.69"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "4.519.99"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'14.58"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "3.913.37"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D18").Value = "'19.67"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "'1.12"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").Value = "68.302.72"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "'433.90"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'14.66"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'3.38"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'87.84"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "'11.44"
$ws.Range("E25").Value = "  +15.09%  "
$ws.Range("D26").Value = "'3.56"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'10.43"
$ws.Range("E27").Value = "  +4.30%  "
$ws.Range("D28").Value = "'37.93"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +4.82%  "
$ws.Range("D30").Value = "'701.85"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'13.38"
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("D32").Value = "'0.130"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").Value = "0.0₃0905"
$ws.Range("E34").Value = "  +33.41%  "
$ws.Range("D35").Value = "'41.50"
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("D36").Value = "'59.57"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").Value = "'5.72"
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("E38").Value = "  -7.14%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "'0.0474"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.74"
$ws.Range("E41").Value = "  +6.95%  "
$ws.Range("D42").Value = "'3.03"
$ws.Range("E42").Value = "  +3.37%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.02"
$ws.Range("E43").Value = "  +9.19%  "
$ws.Range("D44").Value = "'0.341"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "'0.141"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'3.41"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'2.13"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").Value = "'146.10"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").Value = "'3.13"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("E51").Value = "  -2.84%  "
